$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.850.67'
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").Value = '3.410.14'
$ws.Range("E3").Value = '  +0.09%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = 'Normal'
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '411.22'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +0.88%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '130.85'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +2.02%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.620'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -1.67%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.727'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  -0.28%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.134'
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -5.25%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '42.83'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +1.00%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '9.14'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +2.92%  '

$ws.Range("D13").Value = '3.946.43'
$ws.Range("E13").Value = '  -0.36%  '

$ws.Range("E14").Value = '  -0.05%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '0.0000209'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -1.56%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '20.44'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -1.93%  '

$ws.Range("D17").Value = '3.396.69'
$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("E18").Value = '  +2.17%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '12.30'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +1.94%  '

$ws.Range("D20").Value = '61.821.17'
$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '482.56'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +18.27%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '90.13'
$ws.Range("D22").Style = 'Normal'
$ws.Range("E22").Value = '  +0.95%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '3.26'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +2.71%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '13.17'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +0.95%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '3.28'
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +1.57%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '9.71'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +13.47%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '32.99'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +0.71%  '

$ws.Range("E28").Value = '  -0.78%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '8.01'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +5.38%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '2.66'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -2.70%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '11.87'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +0.42%  '

$ws.Range("E32").Value = '  -2.05%  '

$ws.Range("E33").Value = '  -5.15%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '41.08'
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -4.70%  '

$ws.Range("E35").Value = '  -0.68%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '56.77'
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  +5.33%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '0.0488'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  -1.25%  '

$ws.Range("E38").Value = '  +0.01%  '

$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '3.05'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +4.70%  '

$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.331'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +7.14%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '147.73'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +5.37%  '

$ws.Range("E42").Value = '  +0.96%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '3.34'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -0.51%  '

$ws.Range("E44").Value = '  +5.54%  '

$ws.Range("E45").Value = '  +7.56%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '4.26'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +5.43%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '16.55'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +0.06%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '2.34'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +19.67%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '21.95'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  +0.53%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '113.02'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +13.96%  '

$ws.Range("B51").Value = 'Fetch.AI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '1.96'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +20.43%  '
